$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new date column "25-sep" in column BX, following the existing
# "22-sep" column (BW), with the per-family counts for that date.

$ws.Range("BX1").Value = "25-sep"

$ws.Range("BX2").Value = 11
$ws.Range("BX3").Value = 15
$ws.Range("BX4").Value = 11
$ws.Range("BX5").Value = 10
$ws.Range("BX6").Value = 9
$ws.Range("BX7").Value = 13
$ws.Range("BX8").Value = 16
$ws.Range("BX9").Value = 15
$ws.Range("BX10").Value = 15
$ws.Range("BX11").Value = 8

$ws.Range("BX2:BX11").HorizontalAlignment = -4108
$ws.Range("BX2:BX11").NumberFormat = "0"

$ws.Range("BY16").Select()
